$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 15: Morning Glass of Ether
$ws.Range("H15").Value = 596.2778
$ws.Range("I15").Value = 596.2778
$ws.Range("K15").Value = 1788.8334
$ws.Range("M15").Value = -1619.8334
# Row 39: Riches' Brew
$ws.Range("H39").Value = 552.625
$ws.Range("I39").Value = 417.2857
$ws.Range("J39").Value = 1500
$ws.Range("K39").Value = 1251.8571
$ws.Range("L39").Value = 4500
$ws.Range("M39").Value = -955.8571000000002
$ws.Range("N39").Value = -5092
# Row 97: Materia Worth
$ws.Range("H97").Value = 1749.6
$ws.Range("J97").Value = 1749.6
$ws.Range("L97").Value = 5248.799999999999
$ws.Range("N97").Value = -6240.799999999999
# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 6744.769
$ws.Range("J137").Value = 7262.75
$ws.Range("L137").Value = 21788.25
$ws.Range("N137").Value = -26888.25

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
# Row 19: Stadium Envy
$ws.Range("H19").Value = 664132.7
$ws.Range("I19").Value = 1979899
$ws.Range("K19").Value = 1979899
$ws.Range("M19").Value = -1979670
# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 9248.5
$ws.Range("I61").Value = 997.25
$ws.Range("J61").Value = 17499.75
$ws.Range("K61").Value = 997.25
$ws.Range("L61").Value = 17499.75
$ws.Range("M61").Value = -785.25
$ws.Range("N61").Value = -17923.75
# Row 97: Ore for Me
$ws.Range("H97").Value = 646.5
$ws.Range("I97").Value = 649
$ws.Range("J97").Value = 644
$ws.Range("K97").Value = 649
$ws.Range("L97").Value = 644
$ws.Range("M97").Value = -153
$ws.Range("N97").Value = -1636
# Row 102: Smells of Rich Tama-hagane
$ws.Range("H102").Value = 1071
$ws.Range("I102").Value = 916.1667
$ws.Range("K102").Value = 916.1667
$ws.Range("M102").Value = 705.8333
# Row 116: No Scope
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()
# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 8733
$ws.Range("I132").Value = 1333
$ws.Range("K132").Value = 3999
$ws.Range("M132").Value = -1469
# Row 135: Forgiveness for My Shins
$ws.Range("H135").Value = 49857.25
$ws.Range("J135").Value = 49857.25
$ws.Range("L135").Value = 49857.25
$ws.Range("N135").Value = -59997.25
# Row 136: Metal with Mettle
$ws.Range("H136").Value = 9248.5
$ws.Range("I136").Value = 997.25
$ws.Range("J136").Value = 17499.75
$ws.Range("K136").Value = 2991.75
$ws.Range("L136").Value = 52499.25
$ws.Range("M136").Value = -441.75
$ws.Range("N136").Value = -57599.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
# Row 94: High Steal
$ws.Range("H94").Value = 703
$ws.Range("I94").Value = 640.875
$ws.Range("J94").Value = 1200
$ws.Range("K94").Value = 640.875
$ws.Range("L94").Value = 1200
$ws.Range("M94").Value = -189.875
$ws.Range("N94").Value = -2102
# Row 100: And My Axe
$ws.Range("H100").Value = 7407
$ws.Range("J100").Value = 7407
$ws.Range("L100").Value = 7407
$ws.Range("N100").Value = -9571
# Row 105: Ingot to Wing It
$ws.Range("H105").Value = 1811.4
$ws.Range("I105").Value = 1264.25
$ws.Range("K105").Value = 1264.25
$ws.Range("M105").Value = 482.75
# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 5358.3335
$ws.Range("I134").Value = 2977.7778
$ws.Range("K134").Value = 8933.3334
$ws.Range("M134").Value = -6398.3334

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22: Driving Up the Wall
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 1000
$ws.Range("K22").Value = 1000
$ws.Range("M22").Value = -650
# Row 31: Wall Not Found
$ws.Range("H31").Value = 7224.923
$ws.Range("I31").Value = 5301
$ws.Range("J31").Value = 13638
$ws.Range("K31").Value = 5301
$ws.Range("L31").Value = 13638
$ws.Range("M31").Value = -5006
$ws.Range("N31").Value = -14228
# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 7224.923
$ws.Range("I34").Value = 5301
$ws.Range("J34").Value = 13638
$ws.Range("K34").Value = 5301
$ws.Range("L34").Value = 13638
$ws.Range("M34").Value = -5099
$ws.Range("N34").Value = -14042
# Row 99: O Pine
$ws.Range("H99").Value = 6332.6665
$ws.Range("I99").Value = 6999
$ws.Range("K99").Value = 6999
$ws.Range("M99").Value = -5501
# Row 107: Built to Last
$ws.Range("H107").Value = 1195
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 1195
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 1195
$ws.Range("N107").Value = -5035
$ws.Range("M107").ClearContents()
# Row 122: Timber of Tenkonto
$ws.Range("H122").Value = 1874.5
$ws.Range("I122").Value = 1874.5
$ws.Range("K122").Value = 5623.5
$ws.Range("M122").Value = -3173.5
# Row 126: A Better Conductor
$ws.Range("H126").Value = 6332.6665
$ws.Range("I126").Value = 6999
$ws.Range("K126").Value = 20997
$ws.Range("M126").Value = -18527
# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 3363.1052
$ws.Range("I132").Value = 2100.0715
$ws.Range("K132").Value = 6300.2145
$ws.Range("M132").Value = -3770.2145

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 2: Pork Is a Salty Food
$ws.Range("H2").Value = 541.2941
$ws.Range("I2").Value = 275.65384
$ws.Range("K2").Value = 1653.92304
$ws.Range("M2").Value = -1540.92304
# Row 12: Butter Me Up
$ws.Range("H12").Value = 32.11111
$ws.Range("I12").Value = 38.666668
$ws.Range("J12").Value = 28.833334
$ws.Range("K12").Value = 116.000004
$ws.Range("L12").Value = 86.50000199999999
$ws.Range("M12").Value = 56.999996
$ws.Range("N12").Value = -432.500002
# Row 80: Saucy for a Suitor
$ws.Range("H80").Value = 1141
$ws.Range("I80").Value = 1141
$ws.Range("K80").Value = 3423
$ws.Range("M80").Value = -2487
# Row 83: Saved by the Sauce (L)
$ws.Range("H83").Value = 1141
$ws.Range("I83").Value = 1141
$ws.Range("K83").Value = 10269
$ws.Range("M83").Value = -5589
# Row 121: A Cookie for Your Troubles
$ws.Range("H121").Value = 1086
$ws.Range("J121").Value = 1422.4
$ws.Range("L121").Value = 4267.200000000001
$ws.Range("N121").Value = -6887.200000000001
# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 3658.3635
$ws.Range("I131").Value = 3960.2222
$ws.Range("K131").Value = 11880.6666
$ws.Range("M131").Value = -6840.6666

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 11: A Ringing Success
$ws.Range("H11").Value = 4503.25
$ws.Range("J11").Value = 6003.5
$ws.Range("L11").Value = 6003.5
$ws.Range("N11").Value = -6281.5
# Row 15: The Tusk at Hand
$ws.Range("H15").Value = 24750
$ws.Range("J15").Value = 24750
$ws.Range("L15").Value = 24750
$ws.Range("N15").Value = -25326
# Row 43: Get the Green Stuff
$ws.Range("H43").Value = 383.5
$ws.Range("I43").Value = 383.5
$ws.Range("K43").Value = 383.5
$ws.Range("M43").Value = -232.5
# Row 81: The Grander Temple
$ws.Range("H81").Value = 24750
$ws.Range("J81").Value = 24750
$ws.Range("L81").Value = 24750
$ws.Range("N81").Value = -26746
# Row 84: Man with a Dragon Earring (L)
$ws.Range("H84").Value = 24750
$ws.Range("J84").Value = 24750
$ws.Range("L84").Value = 74250
$ws.Range("N84").Value = -84234
# Row 97: If I'd a Koppranickel for Every Time...
$ws.Range("H97").Value = 1024.5
$ws.Range("I97").Value = 549.5
$ws.Range("J97").Value = 1499.5
$ws.Range("K97").Value = 549.5
$ws.Range("L97").Value = 1499.5
$ws.Range("M97").Value = -53.5
$ws.Range("N97").Value = -2491.5
# Row 132: On Board for Lar
$ws.Range("H132").Value = 7565.25
$ws.Range("I132").Value = 5012
$ws.Range("K132").Value = 15036
$ws.Range("M132").Value = -12506
# Row 136: Shiny and Good
$ws.Range("H136").Value = 31500
$ws.Range("J136").Value = 31500
$ws.Range("L136").Value = 94500
$ws.Range("N136").Value = -99600

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban
$ws.Range("H7").Value = 10299.2
$ws.Range("I7").Value = 10299.2
$ws.Range("K7").Value = 10299.2
$ws.Range("M7").Value = -10187.2
# Row 23: Back in the Band
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
# Row 40: Best Served Toad
$ws.Range("H40").Value = 11499
$ws.Range("I40").Value = 11249
$ws.Range("J40").Value = 11999
$ws.Range("K40").Value = 11249
$ws.Range("L40").Value = 11999
$ws.Range("M40").Value = -11113
$ws.Range("N40").Value = -12271
# Row 41: The Hand that Bleeds
$ws.Range("H41").Value = 750
$ws.Range("I41").Value = 750
$ws.Range("K41").Value = 750
$ws.Range("M41").Value = -312
# Row 47: Springtime for Coerthas
$ws.Range("H47").Value = 10000
$ws.Range("I47").Value = 10000
$ws.Range("K47").Value = 10000
$ws.Range("M47").Value = -9510
# Row 52: The Tao of Rabbits
$ws.Range("H52").Value = 10000
$ws.Range("I52").Value = 10000
$ws.Range("K52").Value = 10000
$ws.Range("M52").Value = -9767
# Row 100: Tiger in the Sack
$ws.Range("H100").Value = 7754.2
$ws.Range("I100").Value = 5943.25
$ws.Range("K100").Value = 5943.25
$ws.Range("M100").Value = -5402.25
# Row 126: Battered Books
$ws.Range("H126").Value = 10299.2
$ws.Range("I126").Value = 10299.2
$ws.Range("K126").Value = 30897.6
$ws.Range("M126").Value = -28427.6

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 6: Burn Me Up
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("N6").ClearContents()
# Row 122: Heavy Armoire
$ws.Range("H122").Value = 1865.4546
$ws.Range("I122").Value = 1485
$ws.Range("J122").Value = 1950
$ws.Range("K122").Value = 4455
$ws.Range("L122").Value = 5850
$ws.Range("M122").Value = -2005
$ws.Range("N122").Value = -10750
# Row 132: Comfy Cabins
$ws.Range("H132").Value = 13874.75
$ws.Range("I132").Value = 9749.5
$ws.Range("K132").Value = 29248.5
$ws.Range("M132").Value = -26718.5
